$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.708.64'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '3.529.95'
$ws.Range("E3").Value = '  -2.31%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.71%  '
$ws.Range("D7").Value = '3.517.85'
$ws.Range("E7").Value = '  -2.53%  '
$ws.Range("E8").Value = '  -3.00%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.197'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.642'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.86%  '
$ws.Range("E13").Value = '  -3.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.84%  '
$ws.Range("D15").Value = '4.093.60'
$ws.Range("E15").Value = '  -2.38%  '
$ws.Range("E16").Value = '  -2.87%  '
$ws.Range("D17").Value = '69.679.13'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '3.530.30'
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("E19").Value = '  -3.07%  '
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '535.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.70%  '
$ws.Range("E22").Value = '  -3.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '95.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("E29").Value = '  -4.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("E31").Value = '  -4.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '64.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.49%  '
$ws.Range("E34").Value = '  -4.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '544.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.18'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.02%  '
$ws.Range("E37").Value = '  +3.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '38.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("D40").Value = '0.0₃0757'
$ws.Range("E40").Value = '  -7.36%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.134'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.350.34'
$ws.Range("E42").Value = '  +3.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.55'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0437'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.18%  '
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.37%  '
